$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# The panel data was re-queried from the PanelApp API; refresh the
# per-gene "time_taken" timestamps on the "data" sheet (column F, rows
# 2-27) to the values recorded by that re-run.
$newTimes = @(
  "2021-10-05 14:34:06.381923",
  "2021-10-05 14:34:06.381931",
  "2021-10-05 14:34:06.381934",
  "2021-10-05 14:34:06.381937",
  "2021-10-05 14:34:06.381940",
  "2021-10-05 14:34:06.381943",
  "2021-10-05 14:34:06.381945",
  "2021-10-05 14:34:06.381948",
  "2021-10-05 14:34:06.381951",
  "2021-10-05 14:34:06.381953",
  "2021-10-05 14:34:06.381956",
  "2021-10-05 14:34:06.381958",
  "2021-10-05 14:34:06.381961",
  "2021-10-05 14:34:06.381964",
  "2021-10-05 14:34:06.381966",
  "2021-10-05 14:34:06.381969",
  "2021-10-05 14:34:06.381972",
  "2021-10-05 14:34:06.381974",
  "2021-10-05 14:34:06.381977",
  "2021-10-05 14:34:06.381980",
  "2021-10-05 14:34:06.381982",
  "2021-10-05 14:34:06.381985",
  "2021-10-05 14:34:06.381988",
  "2021-10-05 14:34:06.381990",
  "2021-10-05 14:34:06.381993",
  "2021-10-05 14:34:06.381996"
)
for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $data.Cells.Item($i + 2, 6).Value = $newTimes[$i]
}

# Add the new "metadata" worksheet, describing the data pull itself.
$ws = $wb.Worksheets.Add()
$ws.Name = "metadata"

# Header row (row 1), starting at column B
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Hypertrichosis syndromes"
$ws.Range("C2").Value = 120
$ws.Range("D2").Value = "'0.33"
$ws.Range("E2").Value = "2021-08-20T00:15:20.190722Z"
$ws.Range("F2").Value = "2021-10-05 14:34:06.378254"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/120/?format=json"

# Apply the same bold/bordered/centered style used by the "data" sheet's
# header row (and its A-column index cell) to the new sheet's header row
# and its A2 index cell.
$data.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Finally, move "metadata" so it comes right after "data" in tab order
$ws.Move($null, $wb.Worksheets.Item("data"))
